$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new log rows (172, 173) to the feed_logs sheet.
$ws.Cells.Item(172, 1).Value = 171
$ws.Cells.Item(172, 2).Value = 1
$ws.Cells.Item(172, 3).Value = "2024-06-18 14:11:57"
$ws.Cells.Item(172, 4).Value = 200
$ws.Cells.Item(172, 5).Value = 15

$ws.Cells.Item(173, 1).Value = 172
$ws.Cells.Item(173, 2).Value = 2
$ws.Cells.Item(173, 3).Value = "2024-06-18 14:11:57"
$ws.Cells.Item(173, 4).Value = 200
$ws.Cells.Item(173, 5).Value = 0
